# Apply the "AA10字典表" and "SY17日志表" dictionary-table additions to
# Sheet1, mirroring the existing repeated block layout used throughout the
# sheet (two blank spacer rows with a tinted fill, a header row whose A
# column holds the merged table name, and a series of field-definition
# rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Two spacer rows (137:138), styled like the existing A132:D133 band
#    (tinted fill, centered) and merged into a single block A137:D138.
# ---------------------------------------------------------------------
$ws.Range("A132:D133").Copy()
$ws.Range("A137:D138").PasteSpecial(-4122)
$ws.Range("A137:D138").Merge()

# ---------------------------------------------------------------------
# 2) AA10 dictionary table (rows 139-144)
# ---------------------------------------------------------------------
$ws.Range("A139:A144").Merge()
$ws.Range("A139:A144").HorizontalAlignment = -4108
$ws.Range("A139:A144").VerticalAlignment = -4108

$ws.Range("A139").Value = "AA10字典表"
$ws.Range("B139").Value = "AAA100"
$ws.Range("C139").Value = "Variable characters (30)"
$ws.Range("D139").Value = "代码类别"

$ws.Range("B140").Value = "AAA101"
$ws.Range("C140").Value = "Variable characters (30)"
$ws.Range("D140").Value = "代码总称"

$ws.Range("B141").Value = "AAA102"
$ws.Range("C141").Value = "Variable characters (30)"
$ws.Range("D141").Value = "代码码值"

$ws.Range("B142").Value = "AAA103"
$ws.Range("C142").Value = "Variable characters (30)"
$ws.Range("D142").Value = "代码名"

$ws.Range("B143").Value = "AAA104"
$ws.Range("C143").Value = "Byte "
$ws.Range("D143").Value = "删除标记"

$ws.Range("B144").Value = "AAA105"
$ws.Range("C144").Value = "Variable characters (30)"
$ws.Range("D144").Value = "备用字段"

# ---------------------------------------------------------------------
# 3) Two more spacer rows (145:146), same tinted fill but NOT centered
#    (left as separate, un-merged cells like the rest of the sheet uses
#    for its divider bands).
# ---------------------------------------------------------------------
$ws.Range("A132:D133").Copy()
$ws.Range("A145:D146").PasteSpecial(-4122)
$ws.Range("A145:D146").HorizontalAlignment = 1

# ---------------------------------------------------------------------
# 4) SY17 log table (rows 147-155)
# ---------------------------------------------------------------------
$ws.Range("A147:A153").Merge()
$ws.Range("A147:A153").HorizontalAlignment = -4108
$ws.Range("A147:A153").VerticalAlignment = -4108

$ws.Range("A147").Value = "SY17日志表"
$ws.Range("B147").Value = "CSY170"
$ws.Range("C147").Value = "Integer"
$ws.Range("D147").Value = "日志ID"

$ws.Range("B148").Value = "CSY171"
$ws.Range("C148").Value = "Variable characters (30)"
$ws.Range("D148").Value = "登陆ID"

$ws.Range("B149").Value = "CSY172"
$ws.Range("C149").Value = "Variable characters (30)"
$ws.Range("D149").Value = "访问程序"

$ws.Range("B150").Value = "CSY173"
$ws.Range("C150").Value = "Variable characters (30)"
$ws.Range("D150").Value = "IP地址"

$ws.Range("B151").Value = "CSY174"
$ws.Range("C151").Value = "Variable characters (30)"
$ws.Range("D151").Value = "时间"

$ws.Range("B152").Value = "CSY175"
$ws.Range("C152").Value = "Variable characters (30)"
$ws.Range("D152").Value = "浏览器信息"

$ws.Range("B153").Value = "CSY176"
$ws.Range("C153").Value = "Variable characters (30)"
$ws.Range("D153").Value = "访问系统信息"

$ws.Range("B154").Value = "CSY177"
$ws.Range("C154").Value = "Variable characters (30)"
$ws.Range("D154").Value = "主机名称"

$ws.Range("B155").Value = "CSY178"
$ws.Range("C155").Value = "Variable characters (30)"
$ws.Range("D155").Value = "mac地址"

# ---------------------------------------------------------------------
# 5) Scroll / selection bookkeeping to mirror the author's final view.
# ---------------------------------------------------------------------
$ws.Range("D155").Select()
